$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D edits to stay as plain text (matches original inlineStr cells),
# since many new values are numeric-looking strings (e.g. "1.00", "8.28") that
# Excel would otherwise auto-convert to numbers.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range('D2').Value = '57.146.08'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '3.258.12'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '397.89'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').Value = '108.69'
$ws.Range('E6').Value = '  -2.16%  '
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  +4.22%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').Value = '39.26'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').Value = '3.776.08'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = '8.28'
$ws.Range('E14').Value = '  +1.96%  '
$ws.Range('D15').Value = '18.96'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = '3.253.40'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('E18').Value = '  +3.52%  '
$ws.Range('D19').Value = '57.020.25'
$ws.Range('E19').Value = '  +1.48%  '
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('E21').Value = '  +5.05%  '
$ws.Range('D22').Value = '12.92'
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').Value = '294.24'
$ws.Range('E23').Value = '  -3.59%  '
$ws.Range('D24').Value = '73.91'
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D25').Value = '3.17'
$ws.Range('E25').Value = '  -2.10%  '
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').Value = '7.89'
$ws.Range('E27').Value = '  -3.60%  '
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').Value = '7.49'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '0.168'
$ws.Range('E30').Value = '  -3.20%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('D33').Value = '11.17'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = '39.97'
$ws.Range('E34').Value = '  +9.67%  '
$ws.Range('D35').Value = '0.0494'
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('D37').Value = '51.29'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').Value = '3.46'
$ws.Range('E39').Value = '  -1.84%  '
$ws.Range('D40').Value = '3.01'
$ws.Range('E40').Value = '  -3.86%  '
$ws.Range('D41').Value = '136.88'
$ws.Range('E41').Value = '  +2.71%  '
$ws.Range('E42').Value = '  +1.42%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').Value = '16.74'
$ws.Range('E45').Value = '  -2.77%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '3.88'
$ws.Range('E46').Value = '  -3.59%  '
$ws.Range('D47').Value = '22.38'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('E48').Value = '  +4.28%  '
$ws.Range('D49').Value = '2.147.40'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').Value = '2.46'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('E51').Value = '  -7.04%  '

# Restore default (General) style on column D so only the value/type changed,
# matching the original workbook formatting.
$ws.Range("D2:D50").Style = "Normal"
